$wb = $excel.ActiveWorkbook

# --- Sheet "mesures": add the missing "S" tag in column C for row 44 ---
$wsMesures = $wb.Worksheets.Item("mesures")
$wsMesures.Activate()
$wsMesures.Range("C44").Value = "S"
$wsMesures.Range("C44").Style = $wsMesures.Range("A44").Style
$wsMesures.Range("E12").Select()

# --- Sheet "library_content": insert a publication-date row and bump the version ---
$wsLib = $wb.Worksheets.Item("library_content")
$wsLib.Activate()

# library_version: 2 -> 3
$wsLib.Range("B2").Value = 3

# Insert a new row 4 for "library_publication_date"
$wsLib.Rows.Item(4).Insert()
$wsLib.Range("A4").Value = "library_publication_date"
$wsLib.Range("B4").Value = (Get-Date -Year 2025 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
$wsLib.Range("B4").NumberFormat = "mm-dd-yy"
$wsLib.Range("B4").HorizontalAlignment = -4131

$wsLib.Range("A4").Select()
